# Add two new resource/waste pairs (Fertilizer/FertilizerWaste and
# Farm/FarmWaste) to the Resources sheet, continuing the existing
# R2x/R2xX pattern used by rows 5-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R25 / Fertilizer
$ws.Range("A13").Value = "R25"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Fertilizer"

# R25X / FertilizerWaste
$ws.Range("A14").Value = "R25X"
$ws.Range("B14").Value = -0.25
$ws.Range("C14").Value = "FertilizerWaste"

# R26 / Farm
$ws.Range("A15").Value = "R26"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Farm"

# R26X / FarmWaste
$ws.Range("A16").Value = "R26X"
$ws.Range("B16").Value = -0.25
$ws.Range("C16").Value = "FarmWaste"

# Match the author's final selection position recorded in the sheet view.
$ws.Range("C19").Select()
